# Apply the changes described in the diff for SI2020_Exposiciones.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 9 updates ---
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 8
$ws.Range("K9").Value = 4

# --- Row 17 updates ---
$ws.Range("I17").Value = 8
$ws.Range("L17").Value = 3

# --- Row 23 updates (previously empty) ---
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 9
$ws.Range("I23").Value = 8
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 4
$ws.Range("L23").Value = 3
$ws.Range("M23").Value = 3

# --- Row 24 updates (previously empty) ---
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 9
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = 4
$ws.Range("K24").Value = 4
$ws.Range("L24").Value = 3
$ws.Range("M24").Value = 0

# Recalculate formulas so the dependent N-column totals update
$excel.Calculate()

# --- View state: frozen pane anchor and active selection ---
$ws.Activate()
$panes = $excel.ActiveWindow.Panes
$panes.Item(2).ScrollRow = 8
$panes.Item(2).ScrollColumn = 1

$ws.Range("I24").Select()
